$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price looks like a plain number need to be pre-formatted as
# Text so Excel keeps them as strings (matching the source data) instead of
# silently converting them to numeric values.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
$textRange = $ws.Range($textCells[0])
foreach ($addr in $textCells[1..($textCells.Length - 1)]) {
    $textRange = $excel.Union($textRange, $ws.Range($addr))
}
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '23.704.43'
$ws.Range("E2").Value = '  +1.50%  '
$ws.Range("D3").Value = '1.652.10'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '0.9998'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '304.04'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '0.3815'
$ws.Range("E7").Value = '  +1.77%  '
$ws.Range("D8").Value = '51.45'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '0.3607'
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '1.249'
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").Value = '0.08235'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '22.60'
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").Value = '6.544'
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = '7.405'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '0.00001233'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '1.650.67'
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = '97.11'
$ws.Range("D19").Value = '0.06968'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '6.789'
$ws.Range("E20").Value = '  +4.03%  '
$ws.Range("D21").Value = '17.70'
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '12.62'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '23.709.04'
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("D25").Value = '2.552'
$ws.Range("E25").Value = '  +3.84%  '
$ws.Range("D26").Value = '3.077'
$ws.Range("E26").Value = '  -1.50%  '
$ws.Range("D27").Value = '21.30'
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").Value = '152.28'
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("D29").Value = '5.229'
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("D30").Value = '135.24'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").Value = '1.835.19'
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").Value = '6.883'
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").Value = '1.086'
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("D34").Value = '11.97'
$ws.Range("E34").Value = '  +11.18%  '
$ws.Range("D35").Value = '2.103'
$ws.Range("E35").Value = '  -5.95%  '
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("D37").Value = '0.2520'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("D39").Value = '6.090'
$ws.Range("E39").Value = '  +1.79%  '
$ws.Range("D40").Value = '0.07050'
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").Value = '12.84'
$ws.Range("E41").Value = '  +5.75%  '
$ws.Range("D42").Value = '0.7060'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = '1.340'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '15.94'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").Value = '0.6519'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = '2.340'
$ws.Range("E46").Value = '  +2.45%  '
$ws.Range("D47").Value = '0.9997'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("D49").Value = '0.07988'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '128.11'
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").Value = '  -0.73%  '

# Drop the temporary Text number-format now that the values are stored as
# strings, so the cells end up with no special formatting (same as source).
$textRange.ClearFormats()

Write-Host "Applied cryptos update"